$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 202.16667
$ws.Range("I28").Value = 202.16667
$ws.Range("K28").Value = 202.16667
$ws.Range("M28").Value = 282.83333
$ws.Range("H112").Value = 3174.3794
$ws.Range("J112").Value = 3174.3794
$ws.Range("L112").Value = 9523.138199999999
$ws.Range("N112").Value = -11739.1382
$ws.Range("H116").Value = 5268.4287
$ws.Range("I116").Value = 6131.6665
$ws.Range("J116").Value = 4621
$ws.Range("K116").Value = 6131.6665
$ws.Range("L116").Value = 4621
$ws.Range("M116").Value = -2689.6665
$ws.Range("N116").Value = -11505
$ws.Range("H132").Value = 1002706
$ws.Range("I132").Value = 3070.125
$ws.Range("K132").Value = 9210.375
$ws.Range("M132").Value = -6680.375
$ws.Range("H138").Value = 2549.4546
$ws.Range("J138").Value = 3097.0217
$ws.Range("L138").Value = 9291.0651
$ws.Range("N138").Value = -19571.0651

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1930.6923
$ws.Range("J2").Value = 2624.75
$ws.Range("L2").Value = 2624.75
$ws.Range("N2").Value = -2850.75
$ws.Range("H26").Value = 1656.5
$ws.Range("J26").Value = 1226
$ws.Range("L26").Value = 1226
$ws.Range("N26").Value = -1886
$ws.Range("H61").Value = 3412.647
$ws.Range("I61").Value = 2554.7273
$ws.Range("J61").Value = 4985.5
$ws.Range("K61").Value = 2554.7273
$ws.Range("L61").Value = 4985.5
$ws.Range("M61").Value = -2342.7273
$ws.Range("N61").Value = -5409.5
$ws.Range("H62").Value = 39999.5
$ws.Range("J62").Value = 39999.5
$ws.Range("L62").Value = 39999.5
$ws.Range("N62").Value = -41247.5
$ws.Range("H65").Value = 39999.5
$ws.Range("J65").Value = 39999.5
$ws.Range("L65").Value = 119998.5
$ws.Range("N65").Value = -126238.5
$ws.Range("H116").Value = 1930.6923
$ws.Range("J116").Value = 2624.75
$ws.Range("L116").Value = 2624.75
$ws.Range("N116").Value = -7212.75
$ws.Range("H122").Value = 3123.5334
$ws.Range("I122").Value = 3046.3845
$ws.Range("K122").Value = 9139.1535
$ws.Range("M122").Value = -6689.1535
$ws.Range("H132").Value = 2063
$ws.Range("I132").Value = 1181.3334
$ws.Range("K132").Value = 3544.0002
$ws.Range("M132").Value = -1014.0002
$ws.Range("H136").Value = 3412.647
$ws.Range("I136").Value = 2554.7273
$ws.Range("J136").Value = 4985.5
$ws.Range("K136").Value = 7664.1819
$ws.Range("L136").Value = 14956.5
$ws.Range("M136").Value = -5114.1819
$ws.Range("N136").Value = -20056.5
$ws.Range("H138").Value = 66000
$ws.Range("J138").Value = 66000
$ws.Range("L138").Value = 66000
$ws.Range("N138").Value = -76280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1930.6923
$ws.Range("J3").Value = 2624.75
$ws.Range("L3").Value = 2624.75
$ws.Range("N3").Value = -2852.75
$ws.Range("H62").Value = 60181
$ws.Range("J62").Value = 60181
$ws.Range("L62").Value = 60181
$ws.Range("N62").Value = -61553
$ws.Range("H65").Value = 60181
$ws.Range("J65").Value = 60181
$ws.Range("L65").Value = 180543
$ws.Range("N65").Value = -187407
$ws.Range("H94").Value = 3320.1365
$ws.Range("I94").Value = 2379.2942
$ws.Range("K94").Value = 2379.2942
$ws.Range("M94").Value = -1928.2942
$ws.Range("H134").Value = 2794.926
$ws.Range("I134").Value = 2794.926
$ws.Range("K134").Value = 8384.778
$ws.Range("M134").Value = -5849.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 6423.7627
$ws.Range("J31").Value = 10834.826
$ws.Range("L31").Value = 10834.826
$ws.Range("N31").Value = -11424.826
$ws.Range("H34").Value = 6423.7627
$ws.Range("J34").Value = 10834.826
$ws.Range("L34").Value = 10834.826
$ws.Range("N34").Value = -11238.826
$ws.Range("H86").Value = 9485
$ws.Range("J86").Value = 9412.857
$ws.Range("L86").Value = 9412.857
$ws.Range("N86").Value = -11658.857
$ws.Range("H89").Value = 9485
$ws.Range("J89").Value = 9412.857
$ws.Range("L89").Value = 47064.285
$ws.Range("N89").Value = -58296.285
$ws.Range("H107").Value = 944.5217
$ws.Range("I107").Value = 966.3077
$ws.Range("J107").Value = 916.2
$ws.Range("K107").Value = 966.3077
$ws.Range("L107").Value = 916.2
$ws.Range("M107").Value = 953.6923
$ws.Range("N107").Value = -4756.2
$ws.Range("H134").Value = 1459.5385
$ws.Range("I134").Value = 1536
$ws.Range("K134").Value = 4608
$ws.Range("M134").Value = -2073

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 14599.4
$ws.Range("I50").Value = 998
$ws.Range("K50").Value = 2994
$ws.Range("M50").Value = -2513
$ws.Range("H53").Value = 14599.4
$ws.Range("I53").Value = 998
$ws.Range("K53").Value = 2994
$ws.Range("M53").Value = -2513
$ws.Range("H97").Value = 499
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H113").Value = 778.7
$ws.Range("I113").Value = 465
$ws.Range("K113").Value = 1395
$ws.Range("M113").Value = 775
$ws.Range("H131").Value = 55304.145
$ws.Range("J131").Value = 2801.25
$ws.Range("L131").Value = 8403.75
$ws.Range("N131").Value = -18483.75
$ws.Range("H137").Value = 9874.571
$ws.Range("I137").Value = 18892.666
$ws.Range("K137").Value = 56677.99800000001
$ws.Range("M137").Value = -51577.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1814.8667
$ws.Range("I122").Value = 1390.6666
$ws.Range("J122").Value = 1920.9166
$ws.Range("K122").Value = 4171.9998
$ws.Range("L122").Value = 5762.7498
$ws.Range("M122").Value = -1721.9998
$ws.Range("N122").Value = -10662.7498
$ws.Range("H132").Value = 2416.1853
$ws.Range("I132").Value = 2567.7144
$ws.Range("K132").Value = 7703.1432
$ws.Range("M132").Value = -5173.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 13444.333
$ws.Range("I3").Value = 18166.5
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 18166.5
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = -18054.5
$ws.Range("N3").Value = -4224
$ws.Range("H11").Value = 41666.25
$ws.Range("I11").Value = 33333
$ws.Range("J11").Value = 66666
$ws.Range("K11").Value = 33333
$ws.Range("L11").Value = 66666
$ws.Range("M11").Value = -33193
$ws.Range("N11").Value = -66946
$ws.Range("H15").Value = 13444.333
$ws.Range("I15").Value = 18166.5
$ws.Range("J15").Value = 4000
$ws.Range("K15").Value = 18166.5
$ws.Range("L15").Value = 4000
$ws.Range("M15").Value = -17996.5
$ws.Range("N15").Value = -4340
$ws.Range("H43").Value = 2846067.5
$ws.Range("J43").Value = 5722428.5
$ws.Range("L43").Value = 5722428.5
$ws.Range("N43").Value = -5722814.5
$ws.Range("H74").Value = 37980
$ws.Range("I74").Value = 36966.668
$ws.Range("K74").Value = 36966.668
$ws.Range("M74").Value = -35968.668
$ws.Range("H77").Value = 37980
$ws.Range("I77").Value = 36966.668
$ws.Range("K77").Value = 110900.004
$ws.Range("M77").Value = -105908.004
$ws.Range("H132").Value = 2775.476
$ws.Range("I132").Value = 2948.5
$ws.Range("K132").Value = 8845.5
$ws.Range("M132").Value = -6315.5
$ws.Range("H136").Value = 3863.963
$ws.Range("I136").Value = 3279.8125
$ws.Range("J136").Value = 4713.636
$ws.Range("K136").Value = 9839.4375
$ws.Range("L136").Value = 14140.908
$ws.Range("M136").Value = -7289.4375
$ws.Range("N136").Value = -19240.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 900
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 900
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 392.70587
$ws.Range("I113").Value = 392.70587
$ws.Range("K113").Value = 1178.11761
$ws.Range("M113").Value = 991.88239
$ws.Range("H129").Value = 60000
$ws.Range("J129").Value = 60000
$ws.Range("L129").Value = 60000
$ws.Range("N129").Value = -70000
$ws.Range("H132").Value = 2550.8096
$ws.Range("I132").Value = 2577.2632
$ws.Range("J132").Value = 2299.5
$ws.Range("K132").Value = 7731.7896
$ws.Range("L132").Value = 6898.5
$ws.Range("M132").Value = -5201.7896
$ws.Range("N132").Value = -11958.5
